$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update fixture text with checkmark + final score, and mark Result column
$ws.Range("A2").Value = "FC Midtjylland ✓ - Vejle Boldklub: 5:1"
$ws.Range("G2").Value = "✓"

# Row 3: update fixture text with checkmark + final score, and mark Result column
$ws.Range("A3").Value = "Getafe CF - Real Madrid ✓: 0:1"
$ws.Range("G3").Value = "✓"

# Row 4 (Real CD España) stays unchanged.

# Insert a new row at position 5 for the new fixture, shifting the
# existing rows 5-6 down to 6-7.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "1.FC Slovacko - AC Sparta Prague : 0:0"
$ws.Range("B5").Value = "AC Sparta Prague"
$ws.Range("C5").Value = 60
$ws.Range("D5").Value = 100
$ws.Range("F5").Value = 1.36

# Row 6 (formerly row 5, Red Star Belgrade) - update fixture text and mark Result
$ws.Range("A6").Value = "Red Star Belgrade ✓ - FK IMT Belgrad: 6:1"
$ws.Range("G6").Value = "✓"

# Row 7 (formerly row 6, Fenerbahce) - update fixture text and mark Result
$ws.Range("A7").Value = "Fenerbahce ✓ - Fatih Karagümrük: 2:1"
$ws.Range("G7").Value = "✓"
